$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D4,D5,D6,D9,D10,D14,D19,D20,D21,D22,D23,D24,D26,D30,D33,D36,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D51').NumberFormat = '@'

$ws.Range('D2').Value = '61.565.19'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '3.448.86'
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '581.05'
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('D6').Value = '146.18'
$ws.Range('E6').Value = '  +6.48%  '
$ws.Range('D7').Value = '3.450.38'
$ws.Range('E7').Value = '  +2.01%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '0.476'
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('D10').Value = '7.64'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('E11').Value = '  +2.68%  '
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('D13').Value = '4.038.15'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').Value = '28.01'
$ws.Range('E14').Value = '  +8.88%  '
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('D17').Value = '3.450.80'
$ws.Range('E17').Value = '  +2.10%  '
$ws.Range('D18').Value = '61.657.86'
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('D19').Value = '6.24'
$ws.Range('E19').Value = '  +8.37%  '
$ws.Range('D20').Value = '14.33'
$ws.Range('E20').Value = '  +3.59%  '
$ws.Range('D21').Value = '9.55'
$ws.Range('E21').Value = '  +2.16%  '
$ws.Range('D22').Value = '388.80'
$ws.Range('E22').Value = '  +3.10%  '
$ws.Range('D23').Value = '0.567'
$ws.Range('E23').Value = '  +2.77%  '
$ws.Range('D24').Value = '73.66'
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('D26').Value = '0.996'
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('D28').Value = '3.594.43'
$ws.Range('E28').Value = '  +2.14%  '
$ws.Range('D30').Value = '7.64'
$ws.Range('E30').Value = '  +2.80%  '
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('E32').Value = '  +1.51%  '
$ws.Range('D33').Value = '1.48'
$ws.Range('E33').Value = '  -10.98%  '
$ws.Range('E34').Value = '  +2.04%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = '24.10'
$ws.Range('E36').Value = '  +2.88%  '
$ws.Range('D37').Value = '3.477.20'
$ws.Range('E37').Value = '  +2.12%  '
$ws.Range('D38').Value = '7.02'
$ws.Range('E38').Value = '  +2.75%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').Value = '1.56'
$ws.Range('E39').Value = '  +0.94%  '
$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').Value = '5.14'
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('D41').Value = '166.88'
$ws.Range('E41').Value = '  +1.35%  '
$ws.Range('D42').Value = '0.0785'
$ws.Range('E42').Value = '  +3.22%  '
$ws.Range('D43').Value = '27.49'
$ws.Range('E43').Value = '  +8.12%  '
$ws.Range('D44').Value = '0.806'
$ws.Range('E44').Value = '  +3.76%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = '4.52'
$ws.Range('E45').Value = '  +3.95%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '42.50'
$ws.Range('E46').Value = '  +1.83%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').Value = '1.73'
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').Value = '1.16'
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.570.84'
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('D51').Value = '6.96'
$ws.Range('E51').Value = '  +2.46%  '
